# "Agregar usuarios fuera de sala"
# Admin can add users who are not yet in a room: Server request code 105
# changes from "Salida de admin" to "Agregar usuario a sala", and the
# Client side gains a new notification row (21) plus renamed rows for the
# admin-driven room/user removal notifications.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Server")
$ws2 = $wb.Worksheets.Item("Client")

# --- Server sheet: request 105 repurposed from "Salida de admin" to
#     "Agregar usuario a sala" ----------------------------------------
$ws1.Range("B6").Value = "Agregar usuario a sala"
$ws1.Range("D6").Value = "Usuario"
$ws1.Range("F6").Value = "Sala"

# --- Client sheet: add the new "Usuario agregado por admin" notification
#     row (21) at the bottom of the table -----------------------------
$ws2.Range("A12").Copy()
$ws2.Range("A13").PasteSpecial(-4122)
$ws2.Range("B13").Value = "Usuario agregado por admin"
$ws2.Range("C13").Value = 21
$ws2.Range("E13").Value = "Sala"

# --- Client sheet: rename "Remoción de una sala" -> "Sala eliminada por admin"
$ws2.Range("B5").Value = "Sala eliminada por admin"

# --- Client sheet: row 19 used to announce the removed room ("Sala
#     eliminada"); it now announces the removed user instead, and row 20
#     becomes the explicit notification that a user was removed ---------
$ws2.Range("B11").Value = "Usuario eliminado"
$ws2.Range("E11").Value = "Usuario"
$ws2.Range("B12").Value = "Notificación de usuario eliminado"

# --- Cosmetic: row 5's key column is no longer highlighted -------------
$ws2.Range("A5").Interior.Pattern = -4142

# --- Restore the final selections / active sheet as left by the editor -
$ws1.Range("A6").Select()
$ws2.Range("C12").Select()
